# Insert a new "Industry" column after "Stock Name" (column B) and before
# "Mutual Fund" (previously column C). This shifts the existing C:I columns
# to D:J and populates the new column with industry classifications per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:I to D:J by inserting a new blank column at C.
$ws.Columns.Item(3).Insert()

# Header for the new column, matching the style of the other header cells.
$ws.Cells.Item(1, 3).Value = "Industry"
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)  # xlPasteFormats

# Industry values for each data row (rows 2-37), keyed by row number.
$industries = @{
    2  = "Banks"
    3  = "Banks"
    4  = "Banks"
    5  = "Insurance"
    6  = "Finance"
    7  = "IT - Software"
    8  = "Electrical Equipment"
    9  = "Pharmaceuticals & Biotechnology"
    10 = "Capital Markets"
    11 = "Power"
    12 = "Construction"
    13 = "Diversified FMCG"
    14 = "Ferrous Metals"
    15 = "IT - Services"
    16 = "Telecom - Services"
    17 = "Automobiles"
    18 = "Leisure Services"
    19 = "Realty"
    20 = "Diversified Metals"
    21 = "Finance"
    22 = "Beverages"
    23 = "Pharmaceuticals & Biotechnology"
    24 = "Realty"
    25 = "Realty"
    26 = "Capital Markets"
    27 = "Power"
    28 = "Banks"
    29 = "Power"
    30 = "N.A."
    31 = "Construction"
    32 = "Transport Infrastructure"
    33 = "Insurance"
    34 = "Finance"
    35 = "Banks"
    36 = "Diversified FMCG"
    37 = "Petroleum Products"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
